# Actualización 11 de Mayo - Mañana
# Adds the missing "Rescatables" (make-up exam) student rows to the
# "Rescatables" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$nc       = @(19330051920198, 19330051920199, 18330051920160, 18330051920170, 18330051920177, 18330051920367, 18330051920158, 18330051920395, 18330051920407)
$paterno  = @("HERNANDEZ", "HUERTA", "MARRON", "REYES", "SANDOVAL", "GASCA", "LOPEZ", "GARCIA", "RAMIREZ")
$materno  = @("SILVESTRE", "FLORES", "BLASQUEZ", "MARTINEZ", "GUZMAN", "RUIZ", "SANCHEZ", "TEXCAHUA", "HERNANDEZ")
$nombres  = @("VANESSA", "JOSE GUILLERMO", "DAVID", "SALVADOR", "SAUL BRANDON", "ESTEFANI JHOSSELYNE", "CESAR ADRIAN", "AMISADAI", "YAMILETH")
$largo    = @("CÁLCULO DIFERENCIAL", "CÁLCULO DIFERENCIAL", "MATEMÁTICAS APLICADAS", "MATEMÁTICAS APLICADAS", "MATEMÁTICAS APLICADAS", "MATEMÁTICAS APLICADAS", "MATEMÁTICAS APLICADAS", "MATEMÁTICAS APLICADAS", "MATEMÁTICAS APLICADAS")
$grupo    = @("4BLCM", "4BLCM", "6AEM", "6AEM", "6AEM", "6APM", "6AEM", "6ASM", "6ASM")
$reprob   = @(2, 2, 2, 2, 2, 2, 1, 1, 1)

$count = $nc.Count

for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $nc[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $largo[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $grupo[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $reprob[$i]
}
